# Updated symbol list on Thu Dec 22 07:50:17 UTC 2022 with GitHub Actions
#
# Applies the refreshed coinranking.com scrape to Sheet1:
#  - Price (column D) values drift for (almost) every coin.
#  - Rows 9-16 shift down by one (a new coin, "One", is inserted at rank #8 /
#    row 9) so Coin/Link/Price/Volume for those rows are replaced wholesale.
#  - A couple of "Volume(1h)" (column E) labels lose/gain a "Bestin24h" /
#    "Worstin24h" suffix as the underlying ranking changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while preserving its stored type as TEXT
# (matches the workbook's original inline-string cells, e.g. "5.270" must
# stay the literal text "5.270", not become the number 5.27). Flipping the
# NumberFormat to Text ("@") before the write stops Excel's automatic
# "looks like a number" conversion; re-applying the "Normal" cell style
# afterwards puts formatting back the way it was (General number format,
# original default style) without disturbing the text value just written.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- Column D (Price) drifts for rows with unchanged coin identity ----
Set-TextValue 2  4 "246.31"
Set-TextValue 3  4 "22.87"
Set-TextValue 4  4 "5.267"
Set-TextValue 5  4 "0.05731"
Set-TextValue 6  4 "3.447"
Set-TextValue 7  4 "0.8091"
Set-TextValue 8  4 "0.8855"

Set-TextValue 17 4 "0.006155"
Set-TextValue 18 4 "0.005105"
Set-TextValue 19 4 "0.0009956"

Set-TextValue 21 4 "3.749"
Set-TextValue 22 4 "6.302"
Set-TextValue 23 4 "2.192"
Set-TextValue 24 4 "0.3276"
Set-TextValue 25 4 "0.1319"
Set-TextValue 26 4 "4.142"
Set-TextValue 27 4 "0.0003000"

Set-TextValue 40 4 "0.03901"
Set-TextValue 41 4 "0.006782"
Set-TextValue 42 4 "0.1067"
Set-TextValue 45 4 "0.00005637"
Set-TextValue 47 4 "0.5998"
Set-TextValue 48 4 "0.1734"

# Volume(1h) label change (rank suffix text only; row 47 keeps the same price
# change above, but its "Bestin24h" suffix disappears this run).
Set-TextValue 47 5 "46CoinbaseStockTokenCOIN"

# ---- Rows 9-16: "One" enters the table and bumps everybody else down one ----
# Row 9: WazirX -> One
Set-TextValue 9 2 "One"
Set-TextValue 9 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 9 4 "0.01093"
Set-TextValue 9 5 "8OneONEBestin24h"

# Row 10: MandalaExchangeToken -> WazirX
Set-TextValue 10 2 "WazirX"
Set-TextValue 10 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue 10 4 "0.1445"
Set-TextValue 10 5 "9WazirXWRX"

# Row 11: LiechtensteinCryptoassetsExchange -> MandalaExchangeToken
Set-TextValue 11 2 "MandalaExchangeToken"
Set-TextValue 11 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue 11 4 "0.07376"
Set-TextValue 11 5 "10MandalaExchangeTokenMDX"

# Row 12: BitrueCoin -> LiechtensteinCryptoassetsExchange
Set-TextValue 12 2 "LiechtensteinCryptoassetsExchange"
Set-TextValue 12 3 "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue 12 4 "0.03022"
Set-TextValue 12 5 "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13: BitMartToken -> BitrueCoin
Set-TextValue 13 2 "BitrueCoin"
Set-TextValue 13 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue 13 4 "0.03096"
Set-TextValue 13 5 "12BitrueCoinBTR"

# Row 14: BitForexToken -> BitMartToken
Set-TextValue 14 2 "BitMartToken"
Set-TextValue 14 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue 14 4 "0.09397"
Set-TextValue 14 5 "13BitMartTokenBMX"

# Row 15: CoinExToken -> BitForexToken
Set-TextValue 15 2 "BitForexToken"
Set-TextValue 15 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue 15 4 "0.001587"
Set-TextValue 15 5 "14BitForexTokenBF"

# Row 16: One -> CoinExToken
Set-TextValue 16 2 "CoinExToken"
Set-TextValue 16 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue 16 4 "0.04800"
Set-TextValue 16 5 "15CoinExTokenCET"
